$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the row above into the new row 8
$ws.Range("A7:I7").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)

# Populate the new row's values
$ws.Range("A8").Value = 42654.746481481481
$ws.Range("B8").Value = $true
$ws.Range("C8").Value = 9883.9599999999991
$ws.Range("D8").Value = 9869.16
$ws.Range("E8").Value = 104.43
$ws.Range("F8").Value = 104.74
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = 0.3
$ws.Range("I8").Value = $false
